$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G9").Value = 2.5
$ws.Range("I9").Value = 2.75
$ws.Range("J9").Value = 3.25
$ws.Range("L9").Value = 3.5
$ws.Range("O9").Value = 1.33
$ws.Range("P9").Value = 3.25
$ws.Range("Q9").Value = 2.1
$ws.Range("R9").Value = 1.7
$ws.Range("W9").Value = 7.5
$ws.Range("AA9").Value = 21
$ws.Range("AJ9").Value = 29
$ws.Range("AR9").Value = 67
$ws.Range("AX9").Value = 17
